# Update the GitHub link on the title slide's subtitle placeholder
# from the full URL to the short "@urdans" handle, leaving the
# "Jose Urdaneta" name paragraph untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# "Jose Urdaneta" (13 chars) + paragraph break (1 char) precede the
# URL text, so the URL run starts at character 15.
$urlRange = $tr.Characters(15, 25)
$urlRange.Text = "@urdans"
